$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Style: bold font, thin border all around, centered horizontally, top vertically.
# Build the style once on B1 (keeps the generated cellXfs/borders tables minimal),
# then replicate the exact same style onto A2 via a formats-only paste so both
# cells end up sharing a single new cell format (no stray intermediate xfs).
$b1 = $ws.Range("B1")
$b1.Font.Bold = $true
$b1.HorizontalAlignment = -4108  # xlCenter
$b1.VerticalAlignment = -4160    # xlTop
$b1.Borders.LineStyle = 1        # xlContinuous
$b1.Borders.Weight = 2           # xlThin

$b1.Copy()
$a2 = $ws.Range("A2")
$a2.PasteSpecial(-4122)          # xlPasteFormats
